$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string (e.g. "4.85") that must be
# forced to remain plain text, matching the workbook author's original inline-string
# formatting, instead of being auto-converted to a number by Excel.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D13", "D14", "D15", "D16", "D19", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (price / 1h volume change / occasional name+link swaps).
$ws.Range("D2").Value = '36.289.77'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").Value = '1.922.94'
$ws.Range("E3").Value = '  -2.98%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '238.94'
$ws.Range("E5").Value = '  -2.53%  '
$ws.Range("D6").Value = '0.601'
$ws.Range("E6").Value = '  -4.07%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '55.36'
$ws.Range("E8").Value = '  -6.30%  '
$ws.Range("D9").Value = '0.353'
$ws.Range("E9").Value = '  -5.88%  '
$ws.Range("D10").Value = '0.0818'
$ws.Range("E10").Value = '  -1.47%  '
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("D12").Value = '2.210.93'
$ws.Range("D13").Value = '0.783'
$ws.Range("E13").Value = '  -9.41%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '20.43'
$ws.Range("E14").Value = '  -12.77%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '13.03'
$ws.Range("E15").Value = '  -6.79%  '
$ws.Range("D16").Value = '5.04'
$ws.Range("E16").Value = '  -7.71%  '
$ws.Range("D17").Value = '1.914.77'
$ws.Range("E17").Value = '  -3.70%  '
$ws.Range("D18").Value = '36.176.59'
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").Value = '67.96'
$ws.Range("E19").Value = '  -3.51%  '
$ws.Range("D20").Value = '0.0₃0845'
$ws.Range("E20").Value = '  -3.33%  '
$ws.Range("D21").Value = '223.80'
$ws.Range("E21").Value = '  -4.38%  '
$ws.Range("D22").Value = '4.85'
$ws.Range("E22").Value = '  -8.51%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '2.29'
$ws.Range("E24").Value = '  -10.45%  '
$ws.Range("D26").Value = '8.95'
$ws.Range("E26").Value = '  -9.80%  '
$ws.Range("D27").Value = '159.28'
$ws.Range("E27").Value = '  -2.10%  '
$ws.Range("D28").Value = '0.128'
$ws.Range("E28").Value = '  -3.89%  '
$ws.Range("D29").Value = '18.86'
$ws.Range("E29").Value = '  -4.88%  '
$ws.Range("D30").Value = '0.115'
$ws.Range("E30").Value = '  -3.73%  '
$ws.Range("D31").Value = '1.07'
$ws.Range("E31").Value = '  -9.28%  '
$ws.Range("D32").Value = '4.43'
$ws.Range("E32").Value = '  -9.17%  '
$ws.Range("D33").Value = '0.0608'
$ws.Range("E33").Value = '  -10.23%  '
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").Value = '4.04'
$ws.Range("E35").Value = '  -8.17%  '
$ws.Range("D36").Value = '5.95'
$ws.Range("E36").Value = '  -4.13%  '
$ws.Range("E37").Value = '  -1.47%  '
$ws.Range("D38").Value = '2.08'
$ws.Range("E38").Value = '  -7.64%  '
$ws.Range("D39").Value = '2.88'
$ws.Range("E39").Value = '  -2.98%  '
$ws.Range("D40").Value = '0.0953'
$ws.Range("E40").Value = '  -0.96%  '
$ws.Range("D41").Value = '2.84'
$ws.Range("E41").Value = '  -1.75%  '
$ws.Range("D42").Value = '0.0205'
$ws.Range("E42").Value = '  -3.73%  '
$ws.Range("D43").Value = '1.12'
$ws.Range("E43").Value = '  -8.97%  '
$ws.Range("D44").Value = '15.18'
$ws.Range("E44").Value = '  -6.15%  '
$ws.Range("D45").Value = '1.315.25'
$ws.Range("E45").Value = '  -3.67%  '
$ws.Range("D46").Value = '0.990'
$ws.Range("E46").Value = '  -9.15%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '83.64'
$ws.Range("E47").Value = '  -9.22%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").Value = '6.93'
$ws.Range("E48").Value = '  -6.96%  '
$ws.Range("D49").Value = '2.80'
$ws.Range("E49").Value = '  -0.84%  '
$ws.Range("D50").Value = '2.102.18'
$ws.Range("E50").Value = '  -2.81%  '
$ws.Range("D51").Value = '42.17'
$ws.Range("E51").Value = '  -6.55%  '

# Re-apply the original (unstyled) cell format to the text-forced cells so the only
# change versus the source workbook is the cell content, not its style index.
$ws.Range("C2").Copy()
foreach ($addr in $textCells) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
